# PROVA - re-introdotta differenza minuti vs ore nel calcolo del delta per le LS; alcuni cambi QOL
# Update sheet "Release Date (RD)" with the refreshed id/release_date/tassativita rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Release Date (RD)")

# id, release_date (serial), tassativita  for rows 2..17
$data = @(
    @(253522, 45910.58333333334, 0),
    @(253472, 45912.58333333334, 0),
    @(253392, 45911.58333333334, 0),
    @(253393, 45912.58333333334, 0),
    @(253376, 45910.58333333334, 0),
    @(253295, 45911.58333333334, "X"),
    @(252529, 45910.58333333334, 0),
    @(252274, 45911.58333333334, 0),
    @(253244, 45912.58333333334, 0),
    @(252397, 45911.58333333334, 0),
    @(253549, 45912.58333333334, 0),
    @(253261, 45910.58333333334, 0),
    @(244743, 45910.58333333334, 0),
    @(253436, 45910.58333333334, 0),
    @(252741, 45911.58333333334, 0),
    @(253668, 45910.58333333334, 0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rowValues[0]

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $rowValues[1]
    $bCell.NumberFormat = "yyyy-mm-dd h:mm:ss"

    $ws.Cells.Item($row, 3).Value = $rowValues[2]
}
